$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# Rename the sheet: LoginTests -> flightsSearch
# ------------------------------------------------------------------
$ws.Name = "flightsSearch"

# ------------------------------------------------------------------
# Wipe the old login-test content/format and the stale hyperlinks so
# we can lay down the new flight-search data table from a clean sheet.
# ------------------------------------------------------------------
$ws.Cells.Clear()
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# Header row (row 1): searchFlightsDataExcel / originCity / destinationCity /
# departureDate / returnDate / numOfAdults
# ------------------------------------------------------------------
$ws.Range("A1").Value = "searchFlightsDataExcel"
$ws.Range("A1").Interior.Color = 65535

$ws.Range("B1").Value = "originCity"
$ws.Range("C1").Value = "destinationCity"
$ws.Range("D1").Value = "departureDate"
$ws.Range("B1:D1").Font.Color = 0
$ws.Range("B1:D1").Interior.Color = 0x9BD6C3

$ws.Range("E1").Value = "returnDate"
$ws.Range("F1").Value = "numOfAdults"
$ws.Range("E1:F1").Interior.Color = 0x9BD6C3

# ------------------------------------------------------------------
# Origin-city column (was hyperlinked e-mail addresses, now plain text)
# ------------------------------------------------------------------
$ws.Range("B2").Value = "London"
$ws.Range("B3").Value = "Madrid"
$ws.Range("B4").Value = "Paris"
$ws.Range("B2:B7").Font.Underline = $false
$ws.Range("B2:B7").Font.ThemeColor = 1

# ------------------------------------------------------------------
# Destination-city column
# ------------------------------------------------------------------
$ws.Range("C4").Value = "Chicago"
$ws.Range("C3").Value = "Glasgow"
$ws.Range("C2").Value = "Mexico City"
$ws.Range("C2:C3").Font.Color = 0

# ------------------------------------------------------------------
# Departure / return dates (stored as quote-prefixed text, same as the
# authored workbook) and the adult-count column.
# ------------------------------------------------------------------
$ws.Range("D2").NumberFormat = "m/d/yy"
$ws.Range("D2").Font.Color = 0
$ws.Range("D2").Value = "'01/05/2016"
$ws.Range("E2").NumberFormat = "m/d/yy"
$ws.Range("E2").Value = "'06/05/2016"

$ws.Range("D3").NumberFormat = "m/d/yy"
$ws.Range("D3").Value = "'01/06/2016"
$ws.Range("E3").NumberFormat = "m/d/yy"
$ws.Range("E3").Value = "'07/06/2016"

$ws.Range("D4").NumberFormat = "m/d/yy"
$ws.Range("D4").Value = "'01/07/2016"
$ws.Range("E4").NumberFormat = "m/d/yy"
$ws.Range("E4").Value = "'08/07/2016"

$ws.Range("F2").Value = "'2"
$ws.Range("F3").Value = "'3"
$ws.Range("F4").Value = "'4"

# ------------------------------------------------------------------
# Second copy of the sheet title further down/right (G5), mirroring A1.
# ------------------------------------------------------------------
$ws.Range("G5").Value = "searchFlightsDataExcel"
$ws.Range("G5").Interior.Color = 65535

# ------------------------------------------------------------------
# Rows 18/19 only ever carried number-format-only placeholder cells;
# restore that formatting after the sheet-wide clear.
# ------------------------------------------------------------------
$ws.Range("F18").NumberFormat = "m/d/yy"
$ws.Range("G18").NumberFormat = "m/d/yy"
$ws.Range("F19").NumberFormat = "m/d/yy"
$ws.Range("G19").NumberFormat = "mm/dd/yy;@"

# ------------------------------------------------------------------
# Column widths: best-fit the columns that now hold real content.
# ------------------------------------------------------------------
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null
$ws.Columns.Item(7).AutoFit() | Out-Null

# ------------------------------------------------------------------
# Selection moves from D6 to F11 in the saved view.
# ------------------------------------------------------------------
$ws.Range("F11").Select()
